# Apply the cryptocurrency price/volume updates described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.106.32"
$ws.Range("E2").Value = "  +0.17%  "
# Row 3
$ws.Range("D3").Value = "1.832.77"
$ws.Range("E3").Value = "  -0.37%  "
# Row 4
$ws.Range("E4").Value = "  +0.11%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.23"
$ws.Range("E5").Value = "  -1.69%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6847"
$ws.Range("E6").Value = "  -1.76%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.15%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3016"
$ws.Range("E8").Value = "  -1.25%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07473"
$ws.Range("E9").Value = "  -3.12%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.10"
$ws.Range("E10").Value = "  -1.68%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07661"
$ws.Range("E11").Value = "  -2.09%  "
# Row 12
$ws.Range("D12").Value = "1.842.48"
$ws.Range("E12").Value = "  +0.05%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.067"
$ws.Range("E13").Value = "  -0.97%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6830"
$ws.Range("E14").Value = "  -0.15%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "87.81"
$ws.Range("E15").Value = "  -5.53%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.157"
$ws.Range("E16").Value = "  -6.58%  "
# Row 17
$ws.Range("D17").Value = "29.109.37"
$ws.Range("E17").Value = "  +0.24%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008182"
$ws.Range("E18").Value = "  -1.20%  "
# Row 19
$ws.Range("D19").Value = "2.082.47"
$ws.Range("E19").Value = "  +0.23%  "
# Row 20
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "227.19"
$ws.Range("E20").Value = "  -6.17%  "
# Row 21
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.54"
$ws.Range("E21").Value = "  -1.52%  "
# Row 22
$ws.Range("E22").Value = "  +0.13%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.435"
$ws.Range("E23").Value = "  -0.76%  "
# Row 24
$ws.Range("E24").Value = "  +0.06%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1462"
$ws.Range("E25").Value = "  -3.04%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.94"
$ws.Range("E26").Value = "  +0.63%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.768"
$ws.Range("E27").Value = "  -0.36%  "
# Row 28
$ws.Range("E28").Value = "  -0.53%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.511"
$ws.Range("E29").Value = "  -2.00%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.280"
$ws.Range("E30").Value = "  +1.39%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.155"
$ws.Range("E31").Value = "  -0.28%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.196"
$ws.Range("E32").Value = "  +0.04%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05172"
$ws.Range("E33").Value = "  +1.15%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7687"
$ws.Range("E34").Value = "  -1.48%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.849"
$ws.Range("E35").Value = "  -0.45%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.133"
$ws.Range("E36").Value = "  -1.17%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.674"
$ws.Range("E37").Value = "  -0.75%  "
# Row 38
$ws.Range("D38").Value = "1.310.50"
$ws.Range("E38").Value = "  +0.96%  "
# Row 39
$ws.Range("E39").Value = "  -1.30%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.723"
$ws.Range("E40").Value = "  +0.76%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9326"
$ws.Range("E41").Value = "  -1.63%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.801"
$ws.Range("E42").Value = "  -5.64%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "104.56"
$ws.Range("E43").Value = "  -3.01%  "
# Row 44
$ws.Range("E44").Value = "  +0.20%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.33"
$ws.Range("E45").Value = "  +1.80%  "
# Row 46
$ws.Range("E46").Value = "  +2.81%  "
# Row 47
$ws.Range("D47").Value = "1.984.52"
$ws.Range("E47").Value = "  +0.35%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.603"
$ws.Range("E48").Value = "  -0.84%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5194"
$ws.Range("E49").Value = "  +0.41%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.778"
$ws.Range("E50").Value = "  +1.03%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05930"
$ws.Range("E51").Value = "  +0.99%  "
